$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 242.5  # H28: 167.11111 -> 242.5
$ws.Cells.Item(28, 9).Value = 206.875  # I28: 209.16667 -> 206.875
$ws.Cells.Item(28, 10).Value = 313.75  # J28: 83 -> 313.75
$ws.Cells.Item(28, 11).Value = 206.875  # K28: 209.16667 -> 206.875
$ws.Cells.Item(28, 12).Value = 313.75  # L28: 83 -> 313.75
$ws.Cells.Item(28, 13).Value = 278.125  # M28: 275.83333 -> 278.125
$ws.Cells.Item(28, 14).Value = -1283.75  # N28: -1053 -> -1283.75
$ws.Cells.Item(92, 8).Value = 622.1905  # H92: 633.2857 -> 622.1905
$ws.Cells.Item(92, 9).Value = 529.5263  # I92: 514.95 -> 529.5263
$ws.Cells.Item(92, 10).Value = 1502.5  # J92: 3000 -> 1502.5
$ws.Cells.Item(92, 11).Value = 529.5263  # K92: 514.95 -> 529.5263
$ws.Cells.Item(92, 12).Value = 1502.5  # L92: 3000 -> 1502.5
$ws.Cells.Item(92, 13).Value = 718.4737  # M92: 733.05 -> 718.4737
$ws.Cells.Item(92, 14).Value = -3998.5  # N92: -5496 -> -3998.5
$ws.Cells.Item(94, 8).Value = 2131.5  # H94: 2167.3333 -> 2131.5
$ws.Cells.Item(94, 9).Value = 2131.5  # I94: 2520.8 -> 2131.5
$ws.Cells.Item(94, 10).Value = 0  # J94: 400 -> 0
$ws.Cells.Item(94, 11).Value = 2131.5  # K94: 2520.8 -> 2131.5
$ws.Cells.Item(94, 12).Value = 0  # L94: 400 -> 0
$ws.Cells.Item(94, 13).Value = -1680.5  # M94: -2069.8 -> -1680.5
$ws.Cells.Item(94, 14).ClearContents()  # N94: was -1302
$ws.Cells.Item(98, 8).Value = 14250  # H98: 18000 -> 14250
$ws.Cells.Item(98, 9).Value = 14250  # I98: 18000 -> 14250
$ws.Cells.Item(98, 11).Value = 14250  # K98: 18000 -> 14250
$ws.Cells.Item(98, 13).Value = -12752  # M98: -16502 -> -12752
$ws.Cells.Item(111, 8).Value = 1376.4445  # H111: 1505.0667 -> 1376.4445
$ws.Cells.Item(111, 9).Value = 1677.7273  # I111: 1883.8889 -> 1677.7273
$ws.Cells.Item(111, 10).Value = 903  # J111: 936.8333 -> 903
$ws.Cells.Item(111, 11).Value = 5033.1819  # K111: 5651.6667 -> 5033.1819
$ws.Cells.Item(111, 12).Value = 2709  # L111: 2810.4999 -> 2709
$ws.Cells.Item(111, 13).Value = -1966.1819  # M111: -2584.6667 -> -1966.1819
$ws.Cells.Item(111, 14).Value = -8843  # N111: -8944.499899999999 -> -8843
$ws.Cells.Item(118, 8).Value = 1235.1111  # H118: 1120 -> 1235.1111
$ws.Cells.Item(118, 9).Value = 550  # I118: 473.33334 -> 550
$ws.Cells.Item(118, 10).Value = 1783.2  # J118: 5000 -> 1783.2
$ws.Cells.Item(118, 11).Value = 1650  # K118: 1420.00002 -> 1650
$ws.Cells.Item(118, 12).Value = 5349.6  # L118: 15000 -> 5349.6
$ws.Cells.Item(118, 13).Value = 7  # M118: 236.9999800000001 -> 7
$ws.Cells.Item(118, 14).Value = -8663.6  # N118: -18314 -> -8663.6
$ws.Cells.Item(121, 8).Value = 2000  # H121: 2500 -> 2000
$ws.Cells.Item(121, 10).Value = 2000  # J121: 2500 -> 2000
$ws.Cells.Item(121, 12).Value = 6000  # L121: 7500 -> 6000
$ws.Cells.Item(121, 14).Value = -9494  # N121: -10994 -> -9494
$ws.Cells.Item(122, 8).Value = 14250  # H122: 18000 -> 14250
$ws.Cells.Item(122, 9).Value = 14250  # I122: 18000 -> 14250
$ws.Cells.Item(122, 11).Value = 42750  # K122: 54000 -> 42750
$ws.Cells.Item(122, 13).Value = -40300  # M122: -51550 -> -40300
$ws.Cells.Item(138, 8).Value = 1406.47  # H138: 1384.051 -> 1406.47
$ws.Cells.Item(138, 9).Value = 910.34375  # I138: 896.8823 -> 910.34375
$ws.Cells.Item(138, 10).Value = 1639.9412  # J138: 1642.8594 -> 1639.9412
$ws.Cells.Item(138, 11).Value = 2731.03125  # K138: 2690.6469 -> 2731.03125
$ws.Cells.Item(138, 12).Value = 4919.8236  # L138: 4928.5782 -> 4919.8236
$ws.Cells.Item(138, 13).Value = 2408.96875  # M138: 2449.3531 -> 2408.96875
$ws.Cells.Item(138, 14).Value = -15199.8236  # N138: -15208.5782 -> -15199.8236

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 1086.6923  # H110: 1086.8077 -> 1086.6923
$ws.Cells.Item(110, 9).Value = 634.4  # I110: 634.6 -> 634.4
$ws.Cells.Item(110, 11).Value = 634.4  # K110: 634.6 -> 634.4
$ws.Cells.Item(110, 13).Value = 1410.6  # M110: 1410.4 -> 1410.6
$ws.Cells.Item(139, 8).Value = 32782.625  # H139: 33190.625 -> 32782.625
$ws.Cells.Item(139, 10).Value = 32782.625  # J139: 33190.625 -> 32782.625
$ws.Cells.Item(139, 12).Value = 32782.625  # L139: 33190.625 -> 32782.625
$ws.Cells.Item(139, 14).Value = -43062.625  # N139: -43470.625 -> -43062.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 45455990  # H99: 62501540 -> 45455990
$ws.Cells.Item(99, 9).Value = 50001390  # I99: 71430050 -> 50001390
$ws.Cells.Item(99, 11).Value = 50001390  # K99: 71430050 -> 50001390
$ws.Cells.Item(99, 13).Value = -49999892  # M99: -71428552 -> -49999892
$ws.Cells.Item(107, 8).Value = 2128.25  # H107: 1753.9166 -> 2128.25
$ws.Cells.Item(107, 9).Value = 1750  # I107: 1245.8572 -> 1750
$ws.Cells.Item(107, 10).Value = 2254.3333  # J107: 2465.2 -> 2254.3333
$ws.Cells.Item(107, 11).Value = 1750  # K107: 1245.8572 -> 1750
$ws.Cells.Item(107, 12).Value = 2254.3333  # L107: 2465.2 -> 2254.3333
$ws.Cells.Item(107, 13).Value = 170  # M107: 674.1428000000001 -> 170
$ws.Cells.Item(107, 14).Value = -6094.3333  # N107: -6305.2 -> -6094.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 37037824  # H16: 34483548 -> 37037824
$ws.Cells.Item(16, 9).Value = 40000732  # I16: 40000764 -> 40000732
$ws.Cells.Item(16, 10).Value = 1490  # J16: 945 -> 1490
$ws.Cells.Item(16, 11).Value = 40000732  # K16: 40000764 -> 40000732
$ws.Cells.Item(16, 12).Value = 1490  # L16: 945 -> 1490
$ws.Cells.Item(16, 13).Value = -40000445  # M16: -40000477 -> -40000445
$ws.Cells.Item(16, 14).Value = -2064  # N16: -1519 -> -2064
$ws.Cells.Item(31, 8).Value = 1552.9333  # H31: 1462.0303 -> 1552.9333
$ws.Cells.Item(31, 9).Value = 1369.091  # I31: 1271.16 -> 1369.091
$ws.Cells.Item(31, 11).Value = 1369.091  # K31: 1271.16 -> 1369.091
$ws.Cells.Item(31, 13).Value = -1074.091  # M31: -976.1600000000001 -> -1074.091
$ws.Cells.Item(34, 8).Value = 1552.9333  # H34: 1462.0303 -> 1552.9333
$ws.Cells.Item(34, 9).Value = 1369.091  # I34: 1271.16 -> 1369.091
$ws.Cells.Item(34, 11).Value = 1369.091  # K34: 1271.16 -> 1369.091
$ws.Cells.Item(34, 13).Value = -1167.091  # M34: -1069.16 -> -1167.091
$ws.Cells.Item(113, 8).Value = 37037824  # H113: 34483548 -> 37037824
$ws.Cells.Item(113, 9).Value = 40000732  # I113: 40000764 -> 40000732
$ws.Cells.Item(113, 10).Value = 1490  # J113: 945 -> 1490
$ws.Cells.Item(113, 11).Value = 40000732  # K113: 40000764 -> 40000732
$ws.Cells.Item(113, 12).Value = 1490  # L113: 945 -> 1490
$ws.Cells.Item(113, 13).Value = -39998562  # M113: -39998594 -> -39998562
$ws.Cells.Item(113, 14).Value = -5830  # N113: -5285 -> -5830

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1204.6666  # H5: 1238.4828 -> 1204.6666
$ws.Cells.Item(5, 9).Value = 1211.7241  # I5: 1247 -> 1211.7241
$ws.Cells.Item(5, 11).Value = 3635.1723  # K5: 3741 -> 3635.1723
$ws.Cells.Item(5, 13).Value = -3523.1723  # M5: -3629 -> -3523.1723
$ws.Cells.Item(6, 8).Value = 1099.3334  # H6: 483.5 -> 1099.3334
$ws.Cells.Item(6, 9).Value = 299.66666  # I6: 281.2 -> 299.66666
$ws.Cells.Item(6, 10).Value = 1499.1666  # J6: 1495 -> 1499.1666
$ws.Cells.Item(6, 11).Value = 898.9999799999999  # K6: 843.5999999999999 -> 898.9999799999999
$ws.Cells.Item(6, 12).Value = 4497.4998  # L6: 4485 -> 4497.4998
$ws.Cells.Item(6, 13).Value = -785.9999799999999  # M6: -730.5999999999999 -> -785.9999799999999
$ws.Cells.Item(6, 14).Value = -4723.4998  # N6: -4711 -> -4723.4998
$ws.Cells.Item(64, 8).Value = 4342.2856  # H64: 4523.4443 -> 4342.2856
$ws.Cells.Item(64, 9).Value = 1946  # I64: 1980 -> 1946
$ws.Cells.Item(64, 10).Value = 4741.6665  # J64: 4673.0586 -> 4741.6665
$ws.Cells.Item(64, 11).Value = 5838  # K64: 5940 -> 5838
$ws.Cells.Item(64, 12).Value = 14224.9995  # L64: 14019.1758 -> 14224.9995
$ws.Cells.Item(64, 13).Value = -5568  # M64: -5670 -> -5568
$ws.Cells.Item(64, 14).Value = -14764.9995  # N64: -14559.1758 -> -14764.9995
$ws.Cells.Item(67, 8).Value = 4342.2856  # H67: 4523.4443 -> 4342.2856
$ws.Cells.Item(67, 9).Value = 1946  # I67: 1980 -> 1946
$ws.Cells.Item(67, 10).Value = 4741.6665  # J67: 4673.0586 -> 4741.6665
$ws.Cells.Item(67, 11).Value = 5838  # K67: 5940 -> 5838
$ws.Cells.Item(67, 12).Value = 14224.9995  # L67: 14019.1758 -> 14224.9995
$ws.Cells.Item(67, 13).Value = -4902  # M67: -5004 -> -4902
$ws.Cells.Item(67, 14).Value = -16096.9995  # N67: -15891.1758 -> -16096.9995
$ws.Cells.Item(76, 8).Value = 6441.2856  # H76: 6255.533 -> 6441.2856
$ws.Cells.Item(76, 10).Value = 6166.5386  # J76: 5987.143 -> 6166.5386
$ws.Cells.Item(76, 12).Value = 18499.6158  # L76: 17961.429 -> 18499.6158
$ws.Cells.Item(76, 14).Value = -19265.6158  # N76: -18727.429 -> -19265.6158
$ws.Cells.Item(79, 8).Value = 6441.2856  # H79: 6255.533 -> 6441.2856
$ws.Cells.Item(79, 10).Value = 6166.5386  # J79: 5987.143 -> 6166.5386
$ws.Cells.Item(79, 12).Value = 18499.6158  # L79: 17961.429 -> 18499.6158
$ws.Cells.Item(79, 14).Value = -21151.6158  # N79: -20613.429 -> -21151.6158
$ws.Cells.Item(114, 8).Value = 553.2727  # H114: 561.7 -> 553.2727
$ws.Cells.Item(114, 9).Value = 578  # I114: 507.5 -> 578
$ws.Cells.Item(114, 10).Value = 510  # J114: 643 -> 510
$ws.Cells.Item(114, 11).Value = 1734  # K114: 1522.5 -> 1734
$ws.Cells.Item(114, 12).Value = 1530  # L114: 1929 -> 1530
$ws.Cells.Item(114, 13).Value = 1520  # M114: 1731.5 -> 1520
$ws.Cells.Item(114, 14).Value = -8038  # N114: -8437 -> -8038
$ws.Cells.Item(131, 8).Value = 15874230  # H131: 20409490 -> 15874230
$ws.Cells.Item(131, 10).Value = 1297.7018  # J131: 1452.8837 -> 1297.7018
$ws.Cells.Item(131, 12).Value = 3893.1054  # L131: 4358.6511 -> 3893.1054
$ws.Cells.Item(131, 14).Value = -13973.1054  # N131: -14438.6511 -> -13973.1054
$ws.Cells.Item(135, 8).Value = 1204.6666  # H135: 1238.4828 -> 1204.6666
$ws.Cells.Item(135, 9).Value = 1211.7241  # I135: 1247 -> 1211.7241
$ws.Cells.Item(135, 11).Value = 10905.5169  # K135: 11223 -> 10905.5169
$ws.Cells.Item(135, 13).Value = -8370.516899999999  # M135: -8688 -> -8370.516899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 429.0909  # H97: 503.33334 -> 429.0909
$ws.Cells.Item(97, 9).Value = 480.625  # I97: 512.8570999999999 -> 480.625
$ws.Cells.Item(97, 10).Value = 291.66666  # J97: 470 -> 291.66666
$ws.Cells.Item(97, 11).Value = 480.625  # K97: 512.8570999999999 -> 480.625
$ws.Cells.Item(97, 12).Value = 291.66666  # L97: 470 -> 291.66666
$ws.Cells.Item(97, 13).Value = 15.375  # M97: -16.85709999999995 -> 15.375
$ws.Cells.Item(97, 14).Value = -1283.66666  # N97: -1462 -> -1283.66666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1411.875  # H7: 1666.6666 -> 1411.875
$ws.Cells.Item(7, 9).Value = 1496.6666  # I7: 2500 -> 1496.6666
$ws.Cells.Item(7, 10).Value = 1361  # J7: 1250 -> 1361
$ws.Cells.Item(7, 11).Value = 1496.6666  # K7: 2500 -> 1496.6666
$ws.Cells.Item(7, 12).Value = 1361  # L7: 1250 -> 1361
$ws.Cells.Item(7, 13).Value = -1384.6666  # M7: -2388 -> -1384.6666
$ws.Cells.Item(7, 14).Value = -1585  # N7: -1474 -> -1585
$ws.Cells.Item(16, 8).Value = 540.0417  # H16: 522.4 -> 540.0417
$ws.Cells.Item(16, 9).Value = 533.8  # I16: 513.0952 -> 533.8
$ws.Cells.Item(16, 11).Value = 533.8  # K16: 513.0952 -> 533.8
$ws.Cells.Item(16, 13).Value = -363.8  # M16: -343.0952 -> -363.8
$ws.Cells.Item(40, 8).Value = 2953.5557  # H40: 3317.8 -> 2953.5557
$ws.Cells.Item(40, 9).Value = 2598.6  # I40: 3000 -> 2598.6
$ws.Cells.Item(40, 11).Value = 2598.6  # K40: 3000 -> 2598.6
$ws.Cells.Item(40, 13).Value = -2462.6  # M40: -2864 -> -2462.6
$ws.Cells.Item(54, 8).Value = 0  # H54: 10084 -> 0
$ws.Cells.Item(54, 10).Value = 0  # J54: 10084 -> 0
$ws.Cells.Item(54, 12).Value = 0  # L54: 10084 -> 0
$ws.Cells.Item(54, 14).ClearContents()  # N54: was -11372
$ws.Cells.Item(55, 8).Value = 200.6  # H55: 207 -> 200.6
$ws.Cells.Item(55, 10).Value = 409.83334  # J55: 476 -> 409.83334
$ws.Cells.Item(55, 12).Value = 409.83334  # L55: 476 -> 409.83334
$ws.Cells.Item(55, 14).Value = -755.83334  # N55: -822 -> -755.83334
$ws.Cells.Item(122, 8).Value = 25001636  # H122: 27779512 -> 25001636
$ws.Cells.Item(122, 9).Value = 27779316  # I122: 31251638 -> 27779316
$ws.Cells.Item(122, 11).Value = 83337948  # K122: 93754914 -> 83337948
$ws.Cells.Item(122, 13).Value = -83335498  # M122: -93752464 -> -83335498
$ws.Cells.Item(126, 8).Value = 1411.875  # H126: 1666.6666 -> 1411.875
$ws.Cells.Item(126, 9).Value = 1496.6666  # I126: 2500 -> 1496.6666
$ws.Cells.Item(126, 10).Value = 1361  # J126: 1250 -> 1361
$ws.Cells.Item(126, 11).Value = 4489.9998  # K126: 7500 -> 4489.9998
$ws.Cells.Item(126, 12).Value = 4083  # L126: 3750 -> 4083
$ws.Cells.Item(126, 13).Value = -2019.9998  # M126: -5030 -> -2019.9998
$ws.Cells.Item(126, 14).Value = -9023  # N126: -8690 -> -9023

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 38183788  # H122: 8668792 -> 38183788
$ws.Cells.Item(122, 9).Value = 43450256  # I122: 10002160 -> 43450256
$ws.Cells.Item(122, 11).Value = 130350768  # K122: 30006480 -> 130350768
$ws.Cells.Item(122, 13).Value = -130348318  # M122: -30004030 -> -130348318
$ws.Cells.Item(126, 8).Value = 43479212  # H126: 45455536 -> 43479212
$ws.Cells.Item(126, 9).Value = 52632450  # I126: 58824444 -> 52632450
$ws.Cells.Item(126, 10).Value = 1342.5  # J126: 1262.4 -> 1342.5
$ws.Cells.Item(126, 11).Value = 157897350  # K126: 176473332 -> 157897350
$ws.Cells.Item(126, 12).Value = 4027.5  # L126: 3787.2 -> 4027.5
$ws.Cells.Item(126, 13).Value = -157894880  # M126: -176470862 -> -157894880
$ws.Cells.Item(126, 14).Value = -8967.5  # N126: -8727.200000000001 -> -8967.5
